# Fill in the new Hotcard entity row (row 3) on the HOTCARD sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOTCARD")

# Migration Date / FI Name / Entity ID - format as Text first so
# Excel doesn't auto-convert the date-looking / numeric-looking
# strings into a date serial or a number.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "11/26/2025"

$ws.Range("B3").Value = "First Bank"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "16665"

# Old Platform
$ws.Range("E3").Value = "Legacy Core"

# Service
$ws.Range("G3").Value = "Premium"

# CS Location
$ws.Range("H3").Value = "Dallas, TX"

# Total Card Count
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1000"

# Using OneCall IVR
$ws.Range("J3").Value = "Yes"
